$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '90.676.87'
$ws.Range("E2").Value = '  -0.44%  '
$ws.Range("D3").Value = '3.115.50'
$ws.Range("E3").Value = '  -1.35%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.54'
$ws.Range("E5").Value = '  +9.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '630.47'
$ws.Range("E6").Value = '  +0.74%  '
$ws.Range("E7").Value = '  +1.60%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.361'
$ws.Range("E8").Value = '  -3.31%  '
$ws.Range("D10").Value = '3.113.54'
$ws.Range("E10").Value = '  -0.99%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.714'
$ws.Range("E11").Value = '  -4.30%  '
$ws.Range("E12").Value = '  -1.40%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '36.46'
$ws.Range("E13").Value = '  +4.62%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000245'
$ws.Range("E14").Value = '  -2.21%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.48'
$ws.Range("E15").Value = '  -1.11%  '
$ws.Range("D16").Value = '90.591.73'
$ws.Range("D17").Value = '3.690.05'
$ws.Range("E17").Value = '  -0.80%  '
$ws.Range("D18").Value = '3.155.11'
$ws.Range("E18").Value = '  -0.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.76'
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.12'
$ws.Range("E20").Value = '  -0.86%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000208'
$ws.Range("E21").Value = '  -4.27%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '441.56'
$ws.Range("E22").Value = '  -1.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.54'
$ws.Range("E23").Value = '  +6.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.91'
$ws.Range("E24").Value = '  -0.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.89'
$ws.Range("E25").Value = '  -3.90%  '
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '87.88'
$ws.Range("E26").Value = '  -1.30%  '
$ws.Range("B27").Value = 'Aptos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.40'
$ws.Range("E27").Value = '  +0.21%  '
$ws.Range("D28").Value = '3.306.17'
$ws.Range("E28").Value = '  +0.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.49'
$ws.Range("E30").Value = '  +3.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.159'
$ws.Range("E31").Value = '  -3.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.194'
$ws.Range("E32").Value = '  +25.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.24'
$ws.Range("E33").Value = '  +2.24%  '
$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.889'
$ws.Range("E34").Value = '  -1.84%  '
$ws.Range("B35").Value = 'dogwifhat'
$ws.Range("C35").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.80'
$ws.Range("E35").Value = '  +1.75%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '508.23'
$ws.Range("E36").Value = '  -3.84%  '
$ws.Range("E37").Value = '  +3.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.08'
$ws.Range("E38").Value = '  +0.49%  '
$ws.Range("E39").Value = '  +1.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.28'
$ws.Range("E40").Value = '  -2.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.410'
$ws.Range("E41").Value = '  +0.48%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.17'
$ws.Range("E42").Value = '  -0.41%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0845'
$ws.Range("E44").Value = '  +4.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.24'
$ws.Range("E45").Value = '  +47.00%  '
$ws.Range("E46").Value = '  -2.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '151.40'
$ws.Range("E47").Value = '  +1.66%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.685'
$ws.Range("E48").Value = '  +6.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '45.13'
$ws.Range("E49").Value = '  +1.89%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.34'
$ws.Range("E50").Value = '  -0.49%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.47'
$ws.Range("E51").Value = '  +1.95%  '
